$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

foreach ($r in 2..6) {
    $ws.Cells.Item($r, 3).Value = 45212
}
